# Apply updated vm_pu values (case with 380 kV) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.029988721144156
$ws.Cells.Item(2, 4).Value = 1.038332840094145
$ws.Cells.Item(2, 5).Value = 0.992614727750844
$ws.Cells.Item(2, 6).Value = 1.047191129705649
$ws.Cells.Item(2, 9).Value = 1.036166910946001
$ws.Cells.Item(2, 10).Value = 1.03513244992056
$ws.Cells.Item(2, 11).Value = 1.041121173131158
$ws.Cells.Item(2, 12).Value = 0.9955398523335997
$ws.Cells.Item(2, 13).Value = 1.049954457957876
$ws.Cells.Item(2, 14).Value = 1.036602456730799

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.030905256063485
$ws.Cells.Item(3, 4).Value = 1.039025691889827
$ws.Cells.Item(3, 5).Value = 0.9936372048519299
$ws.Cells.Item(3, 6).Value = 1.048023638693337
$ws.Cells.Item(3, 9).Value = 1.036342632168735
$ws.Cells.Item(3, 10).Value = 1.035690417612634
$ws.Cells.Item(3, 11).Value = 1.041624326250193
$ws.Cells.Item(3, 12).Value = 0.9963617723202687
$ws.Cells.Item(3, 13).Value = 1.050598707605567
$ws.Cells.Item(3, 14).Value = 1.037161216800995

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.031498474504504
$ws.Cells.Item(4, 4).Value = 1.039473835763006
$ws.Cells.Item(4, 5).Value = 0.9942998659930998
$ws.Cells.Item(4, 6).Value = 1.04856246449611
$ws.Cells.Item(4, 9).Value = 1.036454592726619
$ws.Cells.Item(4, 10).Value = 1.036051001899185
$ws.Cells.Item(4, 11).Value = 1.041949047113049
$ws.Cells.Item(4, 12).Value = 0.9968940712668347
$ws.Cells.Item(4, 13).Value = 1.051015069690849
$ws.Cells.Item(4, 14).Value = 1.037522313158593

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.031747900562054
$ws.Cells.Item(5, 4).Value = 1.039662191422356
$ws.Cells.Item(5, 5).Value = 0.994578699834602
$ws.Cells.Item(5, 6).Value = 1.048789018031708
$ws.Cells.Item(5, 9).Value = 1.036501242954986
$ws.Cells.Item(5, 10).Value = 1.036202480831161
$ws.Cells.Item(5, 11).Value = 1.042085354432032
$ws.Cells.Item(5, 12).Value = 0.9971179600053012
$ws.Cells.Item(5, 13).Value = 1.051189984741067
$ws.Cells.Item(5, 14).Value = 1.037674007208027

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.031789782435121
$ws.Cells.Item(6, 4).Value = 1.03969381454054
$ws.Cells.Item(6, 5).Value = 0.994625531979634
$ws.Cells.Item(6, 6).Value = 1.048827059148673
$ws.Cells.Item(6, 9).Value = 1.036509051200429
$ws.Cells.Item(6, 10).Value = 1.036227908288812
$ws.Cells.Item(6, 11).Value = 1.042108228956143
$ws.Cells.Item(6, 12).Value = 0.9971555583673455
$ws.Cells.Item(6, 13).Value = 1.051219346462009
$ws.Cells.Item(6, 14).Value = 1.037699470775584

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.031501807203535
$ws.Cells.Item(7, 4).Value = 1.039476352755451
$ws.Cells.Item(7, 5).Value = 0.994303590798249
$ws.Cells.Item(7, 6).Value = 1.048565491595078
$ws.Cells.Item(7, 9).Value = 1.036455217712625
$ws.Cells.Item(7, 10).Value = 1.036053026404369
$ws.Cells.Item(7, 11).Value = 1.041950869266397
$ws.Cells.Item(7, 12).Value = 0.9968970624462089
$ws.Cells.Item(7, 13).Value = 1.05101740740157
$ws.Cells.Item(7, 14).Value = 1.037524340538807

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.030298434843291
$ws.Cells.Item(8, 4).Value = 1.038567028460723
$ws.Cells.Item(8, 5).Value = 0.9929600610674297
$ws.Cells.Item(8, 6).Value = 1.047472450734564
$ws.Cells.Item(8, 9).Value = 1.036226657166658
$ws.Cells.Item(8, 10).Value = 1.035321111955233
$ws.Cells.Item(8, 11).Value = 1.041291391977562
$ws.Cells.Item(8, 12).Value = 0.9958175282591056
$ws.Cells.Item(8, 13).Value = 1.050172290037803
$ws.Cells.Item(8, 14).Value = 1.0367913866872

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.028179203465877
$ws.Cells.Item(9, 4).Value = 1.036963389410231
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.045547490643674
$ws.Cells.Item(9, 9).Value = 1.035810589396043
$ws.Cells.Item(9, 10).Value = 1.034027921575584
$ws.Cells.Item(9, 11).Value = 1.040122827687241
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.048679234631487
$ws.Cells.Item(9, 14).Value = 1.035496359828881

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.026767293344355
$ws.Cells.Item(10, 4).Value = 1.035893515669231
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.044265028040286
$ws.Cells.Item(10, 9).Value = 1.035524302561888
$ws.Cells.Item(10, 10).Value = 1.033163525586058
$ws.Cells.Item(10, 11).Value = 1.03933949474273
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.04768135424016
$ws.Cells.Item(10, 14).Value = 1.034630736297903

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.026156148939328
$ws.Cells.Item(11, 4).Value = 1.035430081641368
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.043709926973934
$ws.Cells.Item(11, 9).Value = 1.035398233277339
$ws.Cells.Item(11, 10).Value = 1.032788706627772
$ws.Cells.Item(11, 11).Value = 1.038999299265576
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.04724868191338
$ws.Cells.Item(11, 14).Value = 1.034255385053704

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.025929176969178
$ws.Cells.Item(12, 4).Value = 1.035257917266861
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.043503771335023
$ws.Cells.Item(12, 9).Value = 1.035351089931952
$ws.Cells.Item(12, 10).Value = 1.032649403525911
$ws.Cells.Item(12, 11).Value = 1.038872785280943
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.047087881743869
$ws.Cells.Item(12, 14).Value = 1.03411588412546

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.025977861661548
$ws.Cells.Item(13, 4).Value = 1.035294848188036
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.043547990874203
$ws.Cells.Item(13, 9).Value = 1.035361216610217
$ws.Cells.Item(13, 10).Value = 1.0326792880582
$ws.Cells.Item(13, 11).Value = 1.038899929741463
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.047122377815705
$ws.Cells.Item(13, 14).Value = 1.034145811097213

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.026137386660989
$ws.Cells.Item(14, 4).Value = 1.035415850968709
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.043692885390127
$ws.Cells.Item(14, 9).Value = 1.035394342827914
$ws.Cells.Item(14, 10).Value = 1.032777193388952
$ws.Cells.Item(14, 11).Value = 1.038988844638738
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.04723539187649
$ws.Cells.Item(14, 14).Value = 1.034243855464764

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.026235679847179
$ws.Cells.Item(15, 4).Value = 1.035490401592988
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.04378216418251
$ws.Cells.Item(15, 9).Value = 1.035414711184232
$ws.Cells.Item(15, 10).Value = 1.03283750570051
$ws.Cells.Item(15, 11).Value = 1.039043608165672
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.047305012157872
$ws.Cells.Item(15, 14).Value = 1.034304253426723

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.026807857741566
$ws.Cells.Item(16, 4).Value = 1.035924268766631
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.044301872901107
$ws.Cells.Item(16, 9).Value = 1.035532625069344
$ws.Cells.Item(16, 10).Value = 1.033188390004396
$ws.Cells.Item(16, 11).Value = 1.039362051268653
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.047710057091407
$ws.Cells.Item(16, 14).Value = 1.034655636026567

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.027166829902842
$ws.Cells.Item(17, 4).Value = 1.036196377251789
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.044627930729912
$ws.Cells.Item(17, 9).Value = 1.035606026196165
$ws.Cells.Item(17, 10).Value = 1.033408349285079
$ws.Cells.Item(17, 11).Value = 1.03956153333368
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.047963975744955
$ws.Cells.Item(17, 14).Value = 1.034875907674658

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.027376233648799
$ws.Cells.Item(18, 4).Value = 1.036355076870447
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.044818135263669
$ws.Cells.Item(18, 9).Value = 1.035648636665165
$ws.Cells.Item(18, 10).Value = 1.033536596631952
$ws.Cells.Item(18, 11).Value = 1.039677790600802
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.048112025812608
$ws.Cells.Item(18, 14).Value = 1.035004337147474

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.027447638509276
$ws.Cells.Item(19, 4).Value = 1.036409186469865
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.044882993528749
$ws.Cells.Item(19, 9).Value = 1.035663131252583
$ws.Cells.Item(19, 10).Value = 1.033580316953234
$ws.Cells.Item(19, 11).Value = 1.039717414797072
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.048162497456885
$ws.Cells.Item(19, 14).Value = 1.035048119556627

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.027128313371382
$ws.Cells.Item(20, 4).Value = 1.036167184298977
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.044592945677389
$ws.Cells.Item(20, 9).Value = 1.035598171959419
$ws.Cells.Item(20, 10).Value = 1.033384755022881
$ws.Cells.Item(20, 11).Value = 1.039540140846832
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.047936738499668
$ws.Cells.Item(20, 14).Value = 1.034852279905902

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.026090409571214
$ws.Cells.Item(21, 4).Value = 1.035380219326113
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.043650216626638
$ws.Cells.Item(21, 9).Value = 1.035384596689093
$ws.Cells.Item(21, 10).Value = 1.032748364872895
$ws.Cells.Item(21, 11).Value = 1.038962665556315
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.047202114427383
$ws.Cells.Item(21, 14).Value = 1.034214986008907

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.025438037006429
$ws.Cells.Item(22, 4).Value = 1.034885283086771
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.043057680671246
$ws.Cells.Item(22, 9).Value = 1.035248487811321
$ws.Cells.Item(22, 10).Value = 1.032347786535589
$ws.Cells.Item(22, 11).Value = 1.038598715223081
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.046739727839853
$ws.Cells.Item(22, 14).Value = 1.033813838804415

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.025783852900038
$ws.Cells.Item(23, 4).Value = 1.035147670953679
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.043371776186842
$ws.Cells.Item(23, 9).Value = 1.03532081452442
$ws.Cells.Item(23, 10).Value = 1.032560183467543
$ws.Cells.Item(23, 11).Value = 1.038791734166069
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.046984894622815
$ws.Cells.Item(23, 14).Value = 1.034026537364375

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.027145717268069
$ws.Cells.Item(24, 4).Value = 1.03618037538936
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.044608753853447
$ws.Cells.Item(24, 9).Value = 1.035601721578722
$ws.Cells.Item(24, 10).Value = 1.033395416412572
$ws.Cells.Item(24, 11).Value = 1.039549807491371
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.047949046012733
$ws.Cells.Item(24, 14).Value = 1.034862956435988

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.028726919472328
$ws.Cells.Item(25, 4).Value = 1.037378112666169
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.046044996949349
$ws.Cells.Item(25, 9).Value = 1.035919726235269
$ws.Cells.Item(25, 10).Value = 1.034362646576704
$ws.Cells.Item(25, 11).Value = 1.04042569063546
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.049065673414971
$ws.Cells.Item(25, 14).Value = 1.035831560177897
